$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the three hidden "_xlcn.ModelConnection_For_Libro1..." defined
#    names used by the MySQL-for-Excel Data Model tables, appending a "1"
#    suffix (mirrors what Excel does when a second table pull creates a
#    name collision that gets disambiguated).
# ---------------------------------------------------------------------------
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.articulos").Name = "_xlcn.ModelConnection_For_Libro1zigma.articulos1"
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.tipocontenido").Name = "_xlcn.ModelConnection_For_Libro1zigma.tipocontenido1"
$wb.Names.Item("_xlcn.ModelConnection_For_Libro1zigma.usuarios").Name = "_xlcn.ModelConnection_For_Libro1zigma.usuarios1"

# ---------------------------------------------------------------------------
# 2) New small reference tables documenting "rel_tipousuario_menu" next to
#    the existing catalog tables, placed in column G starting at row 30.
# ---------------------------------------------------------------------------
$ws.Range("G30").Value = "rel_tipousuario_menu"
$ws.Range("G31").Value = "idTipoUsuario"
$ws.Range("G32").Value = "idMenu"

$tbl1 = $ws.ListObjects.Add(1, $ws.Range("G30:G32"), 0, 1)
$tbl1.Name = "Tabla1"
$ws.ListObjects.Item("Tabla1").TableStyle = "TableStyleDark3"

$ws.Range("G34").Value = "catMenu"
$ws.Range("G35").Value = "idMenu"
$ws.Range("G36").Value = "nombreCampo"
$ws.Range("G37").Value = "tipoCampo"
$ws.Range("G38").Value = "nombreAMostrar"

$tbl2 = $ws.ListObjects.Add(1, $ws.Range("G34:G38"), 0, 1)
$tbl2.Name = "Tabla2"
$ws.ListObjects.Item("Tabla2").TableStyle = "TableStyleDark3"

# ---------------------------------------------------------------------------
# 3) Widen column G slightly now that it holds the new table headers, and
#    drop the autofit ("best fit") flag since the width is explicit now.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 22

# ---------------------------------------------------------------------------
# 4) Move the visible selection/scroll position to the newly added table.
# ---------------------------------------------------------------------------
$ws.Range("G34:G38").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
